# Populate Company/NAICS/Industry enrichment columns (B, C, D) on the
# "warnCompanies" sheet, then re-select warnCompanies as the active tab
# with the selection left at A13 (matching the author's final state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("warnCompanies")

$ws.Cells.Item(12, 2).Value = "ABM Aviation"
$ws.Cells.Item(12, 3).Value = "48-49"
$ws.Cells.Item(12, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(13, 2).Value = "ABM Aviation"
$ws.Cells.Item(13, 3).Value = "48-49"
$ws.Cells.Item(13, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(14, 2).Value = "ABM Aviation"
$ws.Cells.Item(14, 3).Value = "48-49"
$ws.Cells.Item(14, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(15, 2).Value = "ABM Aviation"
$ws.Cells.Item(15, 3).Value = "48-49"
$ws.Cells.Item(15, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(16, 3).Value = 71
$ws.Cells.Item(16, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(18, 2).Value = "Accenture"
$ws.Cells.Item(18, 3).Value = 51
$ws.Cells.Item(18, 4).Value = "Information"
$ws.Cells.Item(19, 3).Value = 52
$ws.Cells.Item(19, 4).Value = "Finance and Insurance"
$ws.Cells.Item(21, 3).Value = 62
$ws.Cells.Item(21, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(24, 3).Value = 62
$ws.Cells.Item(24, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(25, 3).Value = 62
$ws.Cells.Item(25, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(40, 3).Value = "44-45"
$ws.Cells.Item(40, 4).Value = "Retail Trade"
$ws.Cells.Item(46, 3).Value = "31-33"
$ws.Cells.Item(46, 4).Value = "Manufacturing"
$ws.Cells.Item(47, 3).Value = 61
$ws.Cells.Item(47, 4).Value = "Educational Services"
$ws.Cells.Item(48, 3).Value = "31-33"
$ws.Cells.Item(48, 4).Value = "Manufacturing"
$ws.Cells.Item(49, 3).Value = "31-33"
$ws.Cells.Item(49, 4).Value = "Manufacturing"
$ws.Cells.Item(77, 3).Value = "48-49"
$ws.Cells.Item(77, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(78, 3).Value = "48-49"
$ws.Cells.Item(78, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(87, 3).Value = 61
$ws.Cells.Item(87, 4).Value = "Educational Services"
$ws.Cells.Item(88, 3).Value = "48-49"
$ws.Cells.Item(88, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(95, 3).Value = 81
$ws.Cells.Item(95, 4).Value = "Other Services (except Public Administration)"
$ws.Cells.Item(172, 3).Value = 62
$ws.Cells.Item(172, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(178, 3).Value = "44-45"
$ws.Cells.Item(178, 4).Value = "Retail Trade"
$ws.Cells.Item(189, 3).Value = "44-45"
$ws.Cells.Item(189, 4).Value = "Retail Trade"
$ws.Cells.Item(190, 3).Value = "44-45"
$ws.Cells.Item(190, 4).Value = "Retail Trade"
$ws.Cells.Item(192, 3).Value = 72
$ws.Cells.Item(192, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(248, 3).Value = "31-33"
$ws.Cells.Item(248, 4).Value = "Manufacturing"
$ws.Cells.Item(249, 3).Value = "31-33"
$ws.Cells.Item(249, 4).Value = "Manufacturing"
$ws.Cells.Item(254, 3).Value = "44-45"
$ws.Cells.Item(254, 4).Value = "Retail Trade"
$ws.Cells.Item(255, 3).Value = "44-45"
$ws.Cells.Item(255, 4).Value = "Retail Trade"
$ws.Cells.Item(278, 3).Value = 62
$ws.Cells.Item(278, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(285, 3).Value = 62
$ws.Cells.Item(285, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(286, 3).Value = 62
$ws.Cells.Item(286, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(287, 3).Value = 62
$ws.Cells.Item(287, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(294, 3).Value = "31-33"
$ws.Cells.Item(294, 4).Value = "Manufacturing"
$ws.Cells.Item(295, 3).Value = "31-33"
$ws.Cells.Item(295, 4).Value = "Manufacturing"
$ws.Cells.Item(299, 3).Value = 92
$ws.Cells.Item(299, 4).Value = "Public Administration"
$ws.Cells.Item(303, 3).Value = 62
$ws.Cells.Item(303, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(326, 3).Value = 62
$ws.Cells.Item(326, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(366, 3).Value = 62
$ws.Cells.Item(366, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(407, 3).Value = 51
$ws.Cells.Item(407, 4).Value = "Information"
$ws.Cells.Item(492, 3).Value = 62
$ws.Cells.Item(492, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(559, 3).Value = 62
$ws.Cells.Item(559, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(598, 3).Value = 62
$ws.Cells.Item(598, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(606, 3).Value = 62
$ws.Cells.Item(606, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(628, 3).Value = 62
$ws.Cells.Item(628, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(700, 3).Value = 62
$ws.Cells.Item(700, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(852, 3).Value = 62
$ws.Cells.Item(852, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(882, 3).Value = 62
$ws.Cells.Item(882, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(883, 3).Value = 62
$ws.Cells.Item(883, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(887, 3).Value = "31-33"
$ws.Cells.Item(887, 4).Value = "Manufacturing"
$ws.Cells.Item(970, 3).Value = 62
$ws.Cells.Item(970, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1047, 3).Value = 62
$ws.Cells.Item(1047, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1100, 3).Value = 62
$ws.Cells.Item(1100, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1163, 3).Value = 62
$ws.Cells.Item(1163, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1167, 2).Value = "Shaw Industries"
$ws.Cells.Item(1167, 3).Value = "31-33"
$ws.Cells.Item(1167, 4).Value = "Manufacturing"
$ws.Cells.Item(1168, 2).Value = "Shaw Industries"
$ws.Cells.Item(1168, 3).Value = "31-33"
$ws.Cells.Item(1168, 4).Value = "Manufacturing"
$ws.Cells.Item(1169, 2).Value = "Shaw Industries"
$ws.Cells.Item(1169, 3).Value = "31-33"
$ws.Cells.Item(1169, 4).Value = "Manufacturing"
$ws.Cells.Item(1170, 2).Value = "Shaw Industries"
$ws.Cells.Item(1170, 3).Value = "31-33"
$ws.Cells.Item(1170, 4).Value = "Manufacturing"
$ws.Cells.Item(1171, 2).Value = "Shaw Industries"
$ws.Cells.Item(1171, 3).Value = "31-33"
$ws.Cells.Item(1171, 4).Value = "Manufacturing"
$ws.Cells.Item(1172, 2).Value = "Shaw Industries"
$ws.Cells.Item(1172, 3).Value = "31-33"
$ws.Cells.Item(1172, 4).Value = "Manufacturing"
$ws.Cells.Item(1173, 2).Value = "Shaw Industries"
$ws.Cells.Item(1173, 3).Value = "31-33"
$ws.Cells.Item(1173, 4).Value = "Manufacturing"
$ws.Cells.Item(1174, 2).Value = "Shaw Industries"
$ws.Cells.Item(1174, 3).Value = "31-33"
$ws.Cells.Item(1174, 4).Value = "Manufacturing"
$ws.Cells.Item(1175, 2).Value = "Shaw Industries"
$ws.Cells.Item(1175, 3).Value = "31-33"
$ws.Cells.Item(1175, 4).Value = "Manufacturing"
$ws.Cells.Item(1176, 2).Value = "Shaw Industries"
$ws.Cells.Item(1176, 3).Value = "31-33"
$ws.Cells.Item(1176, 4).Value = "Manufacturing"
$ws.Cells.Item(1177, 2).Value = "Shaw Industries"
$ws.Cells.Item(1177, 3).Value = "31-33"
$ws.Cells.Item(1177, 4).Value = "Manufacturing"
$ws.Cells.Item(1178, 2).Value = "Shaw Industries"
$ws.Cells.Item(1178, 3).Value = "31-33"
$ws.Cells.Item(1178, 4).Value = "Manufacturing"
$ws.Cells.Item(1179, 2).Value = "Shaw Industries"
$ws.Cells.Item(1179, 3).Value = "31-33"
$ws.Cells.Item(1179, 4).Value = "Manufacturing"
$ws.Cells.Item(1198, 3).Value = 62
$ws.Cells.Item(1198, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1201, 2).Value = "Sodexo"
$ws.Cells.Item(1201, 3).Value = 72
$ws.Cells.Item(1201, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1202, 2).Value = "Sodexo"
$ws.Cells.Item(1202, 3).Value = 72
$ws.Cells.Item(1202, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1203, 2).Value = "Sodexo"
$ws.Cells.Item(1203, 3).Value = 72
$ws.Cells.Item(1203, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1204, 2).Value = "Sodexo"
$ws.Cells.Item(1204, 3).Value = 72
$ws.Cells.Item(1204, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1205, 2).Value = "Sodexo"
$ws.Cells.Item(1205, 3).Value = 72
$ws.Cells.Item(1205, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1206, 2).Value = "Sodexo"
$ws.Cells.Item(1206, 3).Value = 72
$ws.Cells.Item(1206, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1207, 2).Value = "Sodexo"
$ws.Cells.Item(1207, 3).Value = 72
$ws.Cells.Item(1207, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1208, 2).Value = "Sodexo"
$ws.Cells.Item(1208, 3).Value = 72
$ws.Cells.Item(1208, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1209, 2).Value = "Sodexo"
$ws.Cells.Item(1209, 3).Value = 72
$ws.Cells.Item(1209, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1210, 2).Value = "Solo Cup"
$ws.Cells.Item(1210, 3).Value = "31-33"
$ws.Cells.Item(1210, 4).Value = "Manufacturing"
$ws.Cells.Item(1213, 3).Value = 71
$ws.Cells.Item(1213, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1214, 3).Value = 71
$ws.Cells.Item(1214, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1216, 3).Value = 62
$ws.Cells.Item(1216, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1217, 3).Value = 62
$ws.Cells.Item(1217, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1218, 3).Value = 62
$ws.Cells.Item(1218, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1219, 3).Value = 62
$ws.Cells.Item(1219, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1228, 3).Value = 62
$ws.Cells.Item(1228, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1278, 3).Value = 52
$ws.Cells.Item(1278, 4).Value = "Finance and Insurance"
$ws.Cells.Item(1279, 3).Value = 72
$ws.Cells.Item(1279, 4).Value = "Accommodation and Food Services"
$ws.Cells.Item(1299, 3).Value = "44-45"
$ws.Cells.Item(1299, 4).Value = "Retail Trade"
$ws.Cells.Item(1309, 3).Value = 71
$ws.Cells.Item(1309, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1315, 3).Value = "31-33"
$ws.Cells.Item(1315, 4).Value = "Manufacturing"
$ws.Cells.Item(1317, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1317, 3).Value = 51
$ws.Cells.Item(1317, 4).Value = "Information"
$ws.Cells.Item(1318, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1318, 3).Value = 51
$ws.Cells.Item(1318, 4).Value = "Information"
$ws.Cells.Item(1319, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1319, 3).Value = 51
$ws.Cells.Item(1319, 4).Value = "Information"
$ws.Cells.Item(1320, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1320, 3).Value = 51
$ws.Cells.Item(1320, 4).Value = "Information"
$ws.Cells.Item(1321, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1321, 3).Value = 51
$ws.Cells.Item(1321, 4).Value = "Information"
$ws.Cells.Item(1322, 2).Value = "The Atlanta Journal Constitution"
$ws.Cells.Item(1322, 3).Value = 51
$ws.Cells.Item(1322, 4).Value = "Information"
$ws.Cells.Item(1326, 3).Value = 71
$ws.Cells.Item(1326, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1327, 3).Value = 71
$ws.Cells.Item(1327, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1328, 3).Value = "44-45"
$ws.Cells.Item(1328, 4).Value = "Retail Trade"
$ws.Cells.Item(1331, 2).Value = "The Home Depot"
$ws.Cells.Item(1331, 3).Value = "44-45"
$ws.Cells.Item(1331, 4).Value = "Retail Trade"
$ws.Cells.Item(1332, 2).Value = "The Home Depot"
$ws.Cells.Item(1332, 3).Value = "44-45"
$ws.Cells.Item(1332, 4).Value = "Retail Trade"
$ws.Cells.Item(1388, 3).Value = 71
$ws.Cells.Item(1388, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1455, 3).Value = "44-45"
$ws.Cells.Item(1455, 4).Value = "Retail Trade"
$ws.Cells.Item(1462, 3).Value = 56
$ws.Cells.Item(1462, 4).Value = "Administrative Support and Waste Management and Remediation Services"
$ws.Cells.Item(1463, 3).Value = 56
$ws.Cells.Item(1463, 4).Value = "Administrative Support and Waste Management and Remediation Services"
$ws.Cells.Item(1467, 3).Value = 62
$ws.Cells.Item(1467, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1469, 3).Value = 71
$ws.Cells.Item(1469, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1481, 3).Value = 62
$ws.Cells.Item(1481, 4).Value = "Health Care and Social Assistance"
$ws.Cells.Item(1506, 3).Value = "48-49"
$ws.Cells.Item(1506, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(1507, 3).Value = "48-49"
$ws.Cells.Item(1507, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(1508, 3).Value = "48-49"
$ws.Cells.Item(1508, 4).Value = "Transportation and Warehousing"
$ws.Cells.Item(1510, 3).Value = "31-33"
$ws.Cells.Item(1510, 4).Value = "Manufacturing"
$ws.Cells.Item(1515, 2).Value = "YogaWorks"
$ws.Cells.Item(1515, 3).Value = 71
$ws.Cells.Item(1515, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1516, 2).Value = "YogaWorks"
$ws.Cells.Item(1516, 3).Value = 71
$ws.Cells.Item(1516, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1517, 2).Value = "YogaWorks"
$ws.Cells.Item(1517, 3).Value = 71
$ws.Cells.Item(1517, 4).Value = "Arts, Entertainment, and Recreation"
$ws.Cells.Item(1518, 2).Value = "YogaWorks"
$ws.Cells.Item(1518, 3).Value = 71
$ws.Cells.Item(1518, 4).Value = "Arts, Entertainment, and Recreation"

# warnCompanies becomes the active/selected sheet (was warnLogs before),
# with cell A13 selected.
$ws.Activate()
$ws.Range("A13").Select()
